# Add season-record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header style (bold font, thin border, centered) from an existing
# header cell onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-59).
$ws.Range("AD2:AD59").Value = 69
$ws.Range("AE2:AE59").Value = 93
$ws.Range("AF2:AF59").Value = 0
